# Update the "users" worksheet:
#  - header row gains a "role" column (replacing the stray "passwordConfirm")
#  - row 2 email corrected to gomriakrem1@gmail.com and given a "role" of admin
#  - row 3 gets a "role" of All
#  - two new users appended (rows 4 and 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlinks so they can be rebuilt cleanly against the
# corrected e-mail addresses (editing .Address in place would just stack a
# second, duplicate hyperlink on top of the old one).
$ws.Hyperlinks.Delete()

# Header row
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "role"

# Row 2
$ws.Range("A2").Value = "gomriakrem1@gmail.com"
$ws.Range("B2").Value = "azerty123"
$ws.Range("C2").Value = "admin"

# Row 3
$ws.Range("A3").Value = "Haroungomri@gmail.com"
$ws.Range("B3").Value = "azerty123"
$ws.Range("C3").Value = "All"

# Row 4 (new)
$ws.Range("A4").Value = "SafaAbid@gmail.com"
$ws.Range("B4").Value = "azerty"

# Row 5 (new)
$ws.Range("A5").Value = "achref.gomri@gmail.com"

# Rebuild the mailto hyperlinks for every e-mail cell
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gomriakrem1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:Haroungomri@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:SafaAbid@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:achref.gomri@gmail.com")

# Hyperlinks.Add doesn't restyle the cell, so re-apply the Hyperlink style
$ws.Range("A2:A5").Style = "Hyperlink"

# Match the author's last selection in the sheet
$null = $ws.Range("C10").Select()
